$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name (row 3, merged C3:D3)
$ws.Range("C3").Value = "Md Apurba Khan"

# Row 7 - __init__ / Attributes are set to input values ...
$ws.Range("E7").Value = "Valid account_number, client_number, balance, date_created, overdraft_limit, and overdraft_rate"
$ws.Range("F7").Value = "ChequingAccount(3001, 1001, 500, date(2023, 1, 1), -200, 0.05)"
$ws.Range("G7").Value = "Instance is created with the given values, and attributes are correctly assigned."

# Row 8 - overdraft limit has invalid type.
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'ChequingAccount(3001, 1001, 500, date(2023, 1, 1), "invalid", 0.05)'
$ws.Range("G8").Value = "Overdraft limit defaults to -100."

# Row 9 - overdraft rate has invalid type.
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'ChequingAccount(3001, 1001, 500, date(2023, 1, 1), -200, "invalid")'
$ws.Range("G9").Value = "Overdraft rate defaults to 0.05."

# Row 10 - date created has invalid type
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'ChequingAccount(3001, 1001, 500, "invalid_date", -200, 0.05)'
$ws.Range("G10").Value = "Raises a TypeError for invalid date format."

# Row 11 - balance greater than overdraft limit
$ws.Range("E11").Value = "self._balance = 500"
$ws.Range("F11").Value = "get_service_charges()"
$ws.Range("G11").Value = "Returns base service charge (e.g., `$0.50)"

# Row 12 - balance less than overdraft limit
$ws.Range("E12").Value = "self._balance = -300"
$ws.Range("F12").Value = "get_service_charges()"
$ws.Range("G12").Value = "Calculates service charge based on overdraft formula (overdraft amount * overdraft rate) + base charge."

# Row 13 - balance equal to overdraft limit
$ws.Range("E13").Value = "self._balance = -200"
$ws.Range("F13").Value = "get_service_charges()"
$ws.Range("G13").Value = "Applies standard overdraft charge without exceeding the limit."

# Row 14 - appropriate value returned based on attribute values. (__str__)
$ws.Range("E14").Value = "Instance initialized with known values"
$ws.Range("F14").Value = "str(chequing_account)"
$ws.Range("G14").Value = "Returns a formatted string representation of the account with overdraft details."

# Restore selection state to match the authored workbook (active cell G14)
$ws.Range("G14").Select()
